$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price (D) cells we are about to rewrite to remain text,
# so numeric-looking strings (e.g. "1.00", "20.40", "72.564.81") are
# not auto-converted to numbers by Excel (matching the source data,
# which stores prices as text). Applied per contiguous block since a
# single comma-separated multi-area Range only honors the first area.
$ws.Range("D2:D9").NumberFormat = "@"
$ws.Range("D11:D18").NumberFormat = "@"
$ws.Range("D20:D23").NumberFormat = "@"
$ws.Range("D25:D33").NumberFormat = "@"
$ws.Range("D35:D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"

$ws.Range("D2").Value = "72.564.81"
$ws.Range("E2").Value = "  +1.10%  "
$ws.Range("D3").Value = "3.976.45"
$ws.Range("E3").Value = "  -0.51%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "586.29"
$ws.Range("E5").Value = "  +8.88%  "
$ws.Range("D6").Value = "152.45"
$ws.Range("E6").Value = "  +2.13%  "
$ws.Range("D7").Value = "0.675"
$ws.Range("E7").Value = "  -2.81%  "
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "0.743"
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("E10").Value = "  -1.68%  "
$ws.Range("D11").Value = "52.84"
$ws.Range("E11").Value = "  +5.28%  "
$ws.Range("D12").Value = "0.0000317"
$ws.Range("E12").Value = "  -1.51%  "
$ws.Range("D13").Value = "10.74"
$ws.Range("E13").Value = "  +0.68%  "
$ws.Range("D14").Value = "4.607.44"
$ws.Range("E14").Value = "  -0.64%  "
$ws.Range("D15").Value = "3.974.85"
$ws.Range("E15").Value = "  -0.68%  "
$ws.Range("D16").Value = "1.27"
$ws.Range("E16").Value = "  +7.81%  "
$ws.Range("D17").Value = "13.96"
$ws.Range("E17").Value = "  -0.37%  "
$ws.Range("D18").Value = "20.40"
$ws.Range("E18").Value = "  +0.12%  "
$ws.Range("E19").Value = "  -0.19%  "
$ws.Range("D20").Value = "72.305.43"
$ws.Range("E20").Value = "  +0.80%  "
$ws.Range("D21").Value = "426.94"
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").Value = "4.67"
$ws.Range("E22").Value = "  +10.68%  "
$ws.Range("D23").Value = "95.48"
$ws.Range("E23").Value = "  -1.48%  "
$ws.Range("E24").Value = "  -1.58%  "
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").Value = "14.19"
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").Value = "4.36"
$ws.Range("E26").Value = "  +18.14%  "
$ws.Range("D27").Value = "11.16"
$ws.Range("E27").Value = "  -1.11%  "
$ws.Range("D28").Value = "10.62"
$ws.Range("E28").Value = "  -0.49%  "
$ws.Range("D29").Value = "5.91"
$ws.Range("E29").Value = "  +1.08%  "
$ws.Range("D30").Value = "36.22"
$ws.Range("E30").Value = "  -1.43%  "
$ws.Range("D31").Value = "7.71"
$ws.Range("E31").Value = "  +5.26%  "
$ws.Range("D32").Value = "49.33"
$ws.Range("E32").Value = "  +1.55%  "
$ws.Range("D33").Value = "13.38"
$ws.Range("E33").Value = "  +0.38%  "
$ws.Range("E34").Value = "  -0.49%  "
$ws.Range("D35").Value = "679.61"
$ws.Range("E35").Value = "  +0.89%  "
$ws.Range("D36").Value = "68.36"
$ws.Range("E36").Value = "  +4.91%  "
$ws.Range("D37").Value = "0.436"
$ws.Range("E37").Value = "  -1.73%  "
$ws.Range("D38").Value = "0.0₃0853"
$ws.Range("E38").Value = "  +4.83%  "
$ws.Range("B39").Value = "Dai"
$ws.Range("C39").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").Value = "  -0.10%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "0.145"
$ws.Range("E40").Value = "  -2.54%  "
$ws.Range("B41").Value = "WEMIXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D41").Value = "3.34"
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("B43").Value = "ThetaToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D43").Value = "3.29"
$ws.Range("E43").Value = "  -3.20%  "
$ws.Range("B44").Value = "THORChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D44").Value = "10.93"
$ws.Range("E44").Value = "  +11.44%  "
$ws.Range("D45").Value = "0.0483"
$ws.Range("E45").Value = "  -0.55%  "
$ws.Range("D46").Value = "2.73"
$ws.Range("E46").Value = "  +3.34%  "
$ws.Range("E47").Value = "  -1.12%  "
$ws.Range("D48").Value = "3.40"
$ws.Range("E48").Value = "  +1.56%  "
$ws.Range("E49").Value = "  +5.47%  "
$ws.Range("E50").Value = "  -0.54%  "
$ws.Range("E51").Value = "  +7.00%  "
